$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at row 167, shifting the old rows 167-186 down to 169-188
$ws.Rows.Item(167).Insert()
$ws.Rows.Item(167).Insert()

# New row 167 (Sandia, Primera, 2023-03-10)
$ws.Cells.Item(167, 1).Value = 11
$ws.Cells.Item(167, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(167, 3).Value = "Bíobío"
$ws.Cells.Item(167, 4).Value = 44995
$ws.Cells.Item(167, 5).Value = 8
$ws.Cells.Item(167, 6).Value = 100112028
$ws.Cells.Item(167, 7).Value = "Sandia"
$ws.Cells.Item(167, 8).Value = "Sin especificar"
$ws.Cells.Item(167, 9).Value = "Primera"
$ws.Cells.Item(167, 10).Value = 500
$ws.Cells.Item(167, 11).Value = 2500
$ws.Cells.Item(167, 12).Value = 2500
$ws.Cells.Item(167, 13).Value = 2500
$ws.Cells.Item(167, 14).Value = "`$/unidad"
$ws.Cells.Item(167, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(167, 16).Value = 2500
$ws.Cells.Item(167, 17).Value = 1
$ws.Cells.Item(167, 18).Value = "Hortaliza"

# New row 168 (Sandia, Segunda, 2023-03-10)
$ws.Cells.Item(168, 1).Value = 11
$ws.Cells.Item(168, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(168, 3).Value = "Bíobío"
$ws.Cells.Item(168, 4).Value = 44995
$ws.Cells.Item(168, 5).Value = 8
$ws.Cells.Item(168, 6).Value = 100112028
$ws.Cells.Item(168, 7).Value = "Sandia"
$ws.Cells.Item(168, 8).Value = "Sin especificar"
$ws.Cells.Item(168, 9).Value = "Segunda"
$ws.Cells.Item(168, 10).Value = 500
$ws.Cells.Item(168, 11).Value = 2000
$ws.Cells.Item(168, 12).Value = 2000
$ws.Cells.Item(168, 13).Value = 2000
$ws.Cells.Item(168, 14).Value = "`$/unidad"
$ws.Cells.Item(168, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(168, 16).Value = 2000
$ws.Cells.Item(168, 17).Value = 1
$ws.Cells.Item(168, 18).Value = "Hortaliza"
